# Auto-generated script applying scheduled price-refresh updates to Typhon_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1553.5714
$ws.Range("I69").Value = 2300
$ws.Range("K69").Value = 6900
$ws.Range("M69").Value = -6026
$ws.Range("H72").Value = 1553.5714
$ws.Range("I72").Value = 2300
$ws.Range("K72").Value = 20700
$ws.Range("M72").Value = -16332
$ws.Range("H74").Value = 12504149
$ws.Range("I74").Value = 3999.6667
$ws.Range("J74").Value = 17861356
$ws.Range("K74").Value = 3999.6667
$ws.Range("L74").Value = 17861356
$ws.Range("M74").Value = -3063.6667
$ws.Range("N74").Value = -17863228
$ws.Range("H77").Value = 12504149
$ws.Range("I77").Value = 3999.6667
$ws.Range("J77").Value = 17861356
$ws.Range("K77").Value = 19998.3335
$ws.Range("L77").Value = 89306780
$ws.Range("M77").Value = -15318.3335
$ws.Range("N77").Value = -89316140
$ws.Range("H92").Value = 125000950
$ws.Range("I92").Value = 142858200
$ws.Range("K92").Value = 142858200
$ws.Range("M92").Value = -142856952
$ws.Range("H100").Value = 2012.2222
$ws.Range("I100").Value = 802.5
$ws.Range("J100").Value = 2357.8572
$ws.Range("K100").Value = 802.5
$ws.Range("L100").Value = 2357.8572
$ws.Range("M100").Value = -261.5
$ws.Range("N100").Value = -3439.8572
$ws.Range("H113").Value = 34486704
$ws.Range("I113").Value = 83336824
$ws.Range("J113").Value = 4267.4116
$ws.Range("K113").Value = 83336824
$ws.Range("L113").Value = 4267.4116
$ws.Range("M113").Value = -83333570
$ws.Range("N113").Value = -10775.4116
$ws.Range("H132").Value = 2214.6956
$ws.Range("I132").Value = 2279
$ws.Range("K132").Value = 6837
$ws.Range("M132").Value = -4307
$ws.Range("H137").Value = 1296.2273
$ws.Range("I137").Value = 1186.0541
$ws.Range("K137").Value = 3558.1623
$ws.Range("M137").Value = -1008.1623
$ws.Range("H138").Value = 1548.3125
$ws.Range("I138").Value = 651.84375
$ws.Range("J138").Value = 3341.25
$ws.Range("K138").Value = 1955.53125
$ws.Range("L138").Value = 10023.75
$ws.Range("M138").Value = 3184.46875
$ws.Range("N138").Value = -20303.75
$ws.Range("H141").Value = 1210.7916
$ws.Range("I141").Value = 1089.5217
$ws.Range("K141").Value = 3268.5651
$ws.Range("M141").Value = 1911.4349

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1261.9131
$ws.Range("I97").Value = 1377.5294
$ws.Range("J97").Value = 934.3333
$ws.Range("K97").Value = 1377.5294
$ws.Range("L97").Value = 934.3333
$ws.Range("M97").Value = -881.5293999999999
$ws.Range("N97").Value = -1926.3333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1364.2858
$ws.Range("I94").Value = 706.4286
$ws.Range("K94").Value = 706.4286
$ws.Range("M94").Value = -255.4286
$ws.Range("H105").Value = 3575232.8
$ws.Range("I105").Value = 4307.5
$ws.Range("J105").Value = 8336466.5
$ws.Range("K105").Value = 4307.5
$ws.Range("L105").Value = 8336466.5
$ws.Range("M105").Value = -2560.5
$ws.Range("N105").Value = -8339960.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 922.1905
$ws.Range("I107").Value = 248.33333
$ws.Range("J107").Value = 1427.5834
$ws.Range("K107").Value = 248.33333
$ws.Range("L107").Value = 1427.5834
$ws.Range("M107").Value = 1671.66667
$ws.Range("N107").Value = -5267.5834
$ws.Range("H132").Value = 15112.743
$ws.Range("I132").Value = 19344.465
$ws.Range("J132").Value = 4341.091
$ws.Range("K132").Value = 58033.395
$ws.Range("L132").Value = 13023.273
$ws.Range("M132").Value = -55503.395
$ws.Range("N132").Value = -18083.273

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1526.5
$ws.Range("I5").Value = 903.5
$ws.Range("K5").Value = 2710.5
$ws.Range("M5").Value = -2598.5
$ws.Range("H113").Value = 630.5
$ws.Range("I113").Value = 618
$ws.Range("J113").Value = 651.3333
$ws.Range("K113").Value = 1854
$ws.Range("L113").Value = 1953.9999
$ws.Range("M113").Value = 316
$ws.Range("N113").Value = -6293.9999
$ws.Range("H122").Value = 515.7143
$ws.Range("I122").Value = 250.23077
$ws.Range("J122").Value = 947.125
$ws.Range("K122").Value = 2252.07693
$ws.Range("L122").Value = 8524.125
$ws.Range("M122").Value = 197.9230699999998
$ws.Range("N122").Value = -13424.125
$ws.Range("H131").Value = 813.89
$ws.Range("I131").Value = 703
$ws.Range("J131").Value = 817.3196
$ws.Range("K131").Value = 2109
$ws.Range("L131").Value = 2451.9588
$ws.Range("M131").Value = 2931
$ws.Range("N131").Value = -12531.9588
$ws.Range("H135").Value = 1526.5
$ws.Range("I135").Value = 903.5
$ws.Range("K135").Value = 8131.5
$ws.Range("M135").Value = -5596.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 3500
$ws.Range("J35").Value = 3500
$ws.Range("L35").Value = 3500
$ws.Range("N35").Value = -4096
$ws.Range("H58").Value = 4872.5
$ws.Range("I58").Value = 2245
$ws.Range("J58").Value = 7500
$ws.Range("K58").Value = 2245
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -1968
$ws.Range("N58").Value = -8054
$ws.Range("H97").Value = 1300
$ws.Range("I97").Value = 1300
$ws.Range("K97").Value = 1300
$ws.Range("M97").Value = -804

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1870.4706
$ws.Range("I68").Value = 1549.75
$ws.Range("J68").Value = 2155.5557
$ws.Range("K68").Value = 1549.75
$ws.Range("L68").Value = 2155.5557
$ws.Range("M68").Value = -800.75
$ws.Range("N68").Value = -3653.5557
$ws.Range("H71").Value = 1870.4706
$ws.Range("I71").Value = 1549.75
$ws.Range("J71").Value = 2155.5557
$ws.Range("K71").Value = 7748.75
$ws.Range("L71").Value = 10777.7785
$ws.Range("M71").Value = -4004.75
$ws.Range("N71").Value = -18265.7785
$ws.Range("H122").Value = 756769.4399999999
$ws.Range("I122").Value = 1636493.5
$ws.Range("J122").Value = 2720.2856
$ws.Range("K122").Value = 4909480.5
$ws.Range("L122").Value = 8160.8568
$ws.Range("M122").Value = -4907030.5
$ws.Range("N122").Value = -13060.8568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3497249.8
$ws.Range("I107").Value = 825.25
$ws.Range("J107").Value = 9091529
$ws.Range("K107").Value = 2475.75
$ws.Range("L107").Value = 27274587
$ws.Range("M107").Value = -555.75
$ws.Range("N107").Value = -27278427
$ws.Range("H136").Value = 35715972
$ws.Range("I136").Value = 47620670
$ws.Range("J136").Value = 1871.1428
$ws.Range("K136").Value = 142862010
$ws.Range("L136").Value = 5613.428400000001
$ws.Range("M136").Value = -142859460
$ws.Range("N136").Value = -10713.4284
